# new SQL implementation in progress; Also testing calculate_rv_precision
#
# Net effect on Sheet1:
#  - Column J (10th col) width: 21 -> 20
#  - Row 3 data fully replaced with a new source row (and J/K/L cleared,
#    since the new row has no mass/lum/radius/spectral-type values)
#  - A new row is inserted at row 10 with fresh data, the old rows 10-13
#    shift down to become rows 11-14, and the old row 14 is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width: J (column 10) from 21 -> 20 ----------------------------
# Excel's ColumnWidth (character units) is offset from the raw OOXML
# <col width> attribute by ~0.83 on this sheet's default font; 19.17 here
# round-trips to an XML width of exactly 20.
$ws.Columns(10).ColumnWidth = 19.17

# --- Row 3: replace with new source row, drop mass/lum/radius/spectype ----
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2043885295912608512"
$ws.Range("B3").Value = 284.2576611101138
$ws.Range("C3").Value = 32.90073932392183
$ws.Range("D3").Value = 5.116928577423096
$ws.Range("E3").Value = 5.428450107574463
$ws.Range("F3").Value = 4.655242443084717
$ws.Range("G3").Value = 66.5742243900887
$ws.Range("H3").Value = 5932
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "2043885295912608512"
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = "2043885295914530944"

# --- Rows 10-14: insert a fresh row at 10 (old 10-13 shift to 11-14), -----
# --- then drop the row that falls off the bottom (old row 14) -------------
$ws.Rows(10).Insert()
$ws.Rows(15).Delete()

# Populate the new row 10 with the new source row's data
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "4590758227640885248"
$ws.Range("B10").Value = 271.7561269586085
$ws.Range("C10").Value = 30.56233941836619
$ws.Range("D10").Value = 4.892927646636963
$ws.Range("E10").Value = 5.211727142333984
$ws.Range("F10").Value = 4.491597652435303
$ws.Range("G10").Value = 63.53576695184549
$ws.Range("H10").Value = 6009
$ws.Range("M10").NumberFormat = "@"
$ws.Range("M10").Value = "4590758227640885248"
$ws.Range("N10").NumberFormat = "@"
$ws.Range("N10").Value = "4590758227637479040"

# Row 12 (old row 11 shifted down) gains an L value ("K") that the old row
# 12 it displaced didn't have.
$ws.Range("L12").Value = "K"
